$d = $word.ActiveDocument

# The FAQ section ("How much SmartRewards can I get each month?" ... through
# "...every second block 1000 addresses will get paid.") is being removed,
# leaving only the "SMARTREWARDS FAQ'S" heading and the trailing empty
# paragraph in place.

# Locate the start of the text to remove: first FAQ question paragraph.
$startRange = $d.Content
$startRange.Find.Execute("How much SmartRewards can I get each month?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $startRange.Paragraphs(1).Range.Start

# Locate the end of the text to remove: the last FAQ answer paragraph.
$endRange = $d.Content
$endRange.Find.Execute("will get paid.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $endRange.Paragraphs(1).Range.End

$delRange = $d.Range($startPos, $endPos)
$delRange.Delete()
